# Generate Report for Handback
# - Overview sheet: "Ready for handoff" status text becomes "Handback transform failed"
#   (shared string is reused by E3/F3 so editing the text updates both cells)
# - zh-cn / de-de sheets: new "Error Detail" column (P) gets a handback/handoff
#   filename-mismatch message, and that column is widened to fit the text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Update the shared status text - this string is reused by the Overview
# sheet's zh-cn/de-de status columns (E3, F3) as well as the per-locale
# sheets' Status column (C3), so every occurrence must be updated together.
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Widen the "Error Detail" column (P) on both locale sheets to fit the message.
# (ColumnWidth here is in "characters"; the engine adds ~5/6 of a character's
# worth of cell padding when it stores the column's OOXML width, so back that
# off here to land on a stored width of exactly 40.)
$zhcn.Columns.Item(16).ColumnWidth = 39.16666666666666
$dede.Columns.Item(16).ColumnWidth = 39.16666666666666

# Populate the per-locale "Error Detail" cell with the handback/handoff mismatch message.
$zhcn.Range("P3").Value = "Handback file name: auacfqek.yg0 is different with handoff file name: 3eac4320-fa9b-44d7-a317-c4a6629b3ff8.fcd715134a3dfecb88aa83083497bfac487b42e9.zh-cn."
$dede.Range("P3").Value = "Handback file name: auacfqek.yg0 is different with handoff file name: 3eac4320-fa9b-44d7-a317-c4a6629b3ff8.fcd715134a3dfecb88aa83083497bfac487b42e9.de-de."
